$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F3").Value = "Disguise Kit=*Artisan’s Tools"
$ws.Range("G3").Value = "*Gaming Sets"
$ws.Range("H3").Value = "Vehicles (Land)"
$ws.Range("L3").Value = "*Artisan’s Tools"
$ws.Range("N3").Value = "*Artisan’s Tools"
$ws.Range("P3").Value = "*Gaming Sets=Thieves’ Tools"
$ws.Range("Q3").Value = "Vehicles (Water)"
$ws.Range("T3").Value = "*Artisan’s Tools"
$ws.Range("V3").Value = "Disguise Kit=*Musical Instruments"
$ws.Range("W3").Value = "*Musical Instruments"
$ws.Range("X3").Value = "*Artisan’s Tools=Vehicles (Land)"
$ws.Range("Y3").Value = "Thieves’ tools=*Musical Instruments"
$ws.Range("Z3").Value = "Poisoner’s kit"
$ws.Range("AB3").Value = "*Gaming Sets"
$ws.Range("AD3").Value = "*Artisan’s Tools"
$ws.Range("AE3").Value = "*Gaming Sets=Vehicles (Water)"
$ws.Range("AK3").Value = "*Musical Instruments"
$ws.Range("AL3").Value = "*Gaming Sets=Vehicles (Land)"
$ws.Range("AM3").Value = "Thieves’ Tools=*Artisan’s Tools"
$ws.Range("AN3").Value = "*Gaming Sets=Vehicles (Land)"
$ws.Range("AO3").Value = "*Artisan’s Tools"
$ws.Range("AP3").Value = "*Gaming Sets"
$ws.Range("AR3").Value = "*Gaming Sets=Thieves’ Tools"
$ws.Range("F5").Value = "Disguise Kit=Common Clothes=A Tears of Virulence Emblem=A Writ of Free Agency Signed by the Lord Regent=*Artisan’s Tools=A Pouch with 15 gp (payment for services rendered)"
$ws.Range("L5").Value = "*Artisan’s Tools=A Maker’s Mark Chisel used to mark your Handiwork with the Symbol of the Clan of Crafters you learned your Skill from=A Set of Traveler’s Clothes=A Pouch containing 5 gp=A Gem worth 10 gp"
$ws.Range("N5").Value = "A Two-Person Tent=*Artisan’s Tools=A Holy Symbol=A Set of Traveler’s Clothes=A Belt Pouch containing 5 gp"
$ws.Range("T5").Value = "*Artisan’s Tools=Merchant’s Scale=A Set of Fine Clothes=A Belt Pouch containing 10 gp"
$ws.Range("U5").Value = "A Disguise Kit=A Costume=A Pouch containing 10gp"
$ws.Range("V5").Value = "*Musical Instruments=The Favor of an Admirer (Trinket)=A Costume=A Belt Pouch containing 15 gp."
$ws.Range("W5").Value = "One Set of Traveler’s Clothes=*Musical Instruments=Poorly Wrought Maps from your Homeland that Depict where you are in Faerûn=A Small Piece of Jewelry worth 10 gp in the Style of your Homeland’s Craftsmanship=A Pouch containing 5 gp"
$ws.Range("X5").Value = "*Artisan’s Tools=A Shovel=An Iron Pot=A Set of Common Clothes=A Pouch containing 10 gp"
$ws.Range("Y5").Value = "A Battered Alms Box=*Musical Instruments=A Cast-Off Military Jacket=A Set of Common Clothes=A Belt Pouch containing 10 gp"
$ws.Range("AB5").Value = "*Gaming Sets=A Lucky Charm=A Set of Fine Clothes=A Belt Pouch containing 15 gp"
$ws.Range("AD5").Value = "*Artisan’s Tools=A Letter of Introduction from your Guild=A Set of Traveler’s Clothes=A Belt Pouch containing 15 gp"
$ws.Range("AE5").Value = "Fishing tackle=*Gaming Sets=A Set of Common Clothes=Rowboat=A Belt containing 5 gp"
$ws.Range("AK5").Value = "Your Inheritance=A Set of Traveler’s Clothes=*Musical Instruments=A Pouch containing 15 gp"
$ws.Range("AL5").Value = "A Simple Puzzle Box=A Scroll containing the Basic Teachings of the Five Gods=*Gaming Sets=A Set of Common Clothes=A Pouch containing 15 gp"
$ws.Range("AO5").Value = "An Izzet Insignia=*Artisan’s Tools=The Charred and Twisted remains of a Failed Experiment=A Hammer=A Block and Tackle=A Set of Common Clothes=A Belt Pouch containing 5 gp (Azorius 1-Zino Coins)"
